$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 21
$ws.Range("H21").Value = 31666.666
$ws.Range("I21").Value = 31666.666
$ws.Range("K21").Value = 31666.666
$ws.Range("M21").Value = -31198.666
# Row 23
$ws.Range("H23").Value = 31666.666
$ws.Range("I23").Value = 31666.666
$ws.Range("K23").Value = 31666.666
$ws.Range("M23").Value = -31432.666
# Row 58
$ws.Range("H58").Value = 1449.7894
$ws.Range("I58").Value = 146.4
$ws.Range("J58").Value = 2898
$ws.Range("K58").Value = 439.2
$ws.Range("L58").Value = 8694
$ws.Range("M58").Value = -289.2
$ws.Range("N58").Value = -8994
# Row 64
$ws.Range("H64").Value = 3089.103
$ws.Range("I64").Value = 2902.7673
$ws.Range("J64").Value = 3409.6
$ws.Range("K64").Value = 2902.7673
$ws.Range("L64").Value = 3409.6
$ws.Range("M64").Value = -2654.7673
$ws.Range("N64").Value = -3905.6
# Row 67
$ws.Range("H67").Value = 3089.103
$ws.Range("I67").Value = 2902.7673
$ws.Range("J67").Value = 3409.6
$ws.Range("K67").Value = 2902.7673
$ws.Range("L67").Value = 3409.6
$ws.Range("M67").Value = -2044.7673
$ws.Range("N67").Value = -5125.6
# Row 74
$ws.Range("H74").Value = 68851.57000000001
$ws.Range("I74").Value = 104102.445
$ws.Range("K74").Value = 104102.445
$ws.Range("M74").Value = -103166.445
# Row 76
$ws.Range("H76").Value = 7129
$ws.Range("I76").Value = 7129
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 7129
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -6814
$ws.Range("N76").ClearContents()
# Row 77
$ws.Range("H77").Value = 68851.57000000001
$ws.Range("I77").Value = 104102.445
$ws.Range("K77").Value = 520512.225
$ws.Range("M77").Value = -515832.225
# Row 79
$ws.Range("H79").Value = 7129
$ws.Range("I79").Value = 7129
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 7129
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -6037
$ws.Range("N79").ClearContents()
# Row 82
$ws.Range("H82").Value = 8437.571
$ws.Range("I82").Value = 1093.8334
$ws.Range("J82").Value = 52500
$ws.Range("K82").Value = 3281.5002
$ws.Range("L82").Value = 157500
$ws.Range("M82").Value = -2875.5002
$ws.Range("N82").Value = -158312
# Row 85
$ws.Range("H85").Value = 8437.571
$ws.Range("I85").Value = 1093.8334
$ws.Range("J85").Value = 52500
$ws.Range("K85").Value = 3281.5002
$ws.Range("L85").Value = 157500
$ws.Range("M85").Value = -1877.5002
$ws.Range("N85").Value = -160308
# Row 123
$ws.Range("H123").Value = 69111.11
$ws.Range("J123").Value = 69111.11
$ws.Range("L123").Value = 69111.11
$ws.Range("N123").Value = -78911.11
# Row 126
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
# Row 128
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
# Row 130
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
# Row 134
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 123
$ws.Range("H123").Value = 30425.834
$ws.Range("J123").Value = 30425.834
$ws.Range("L123").Value = 30425.834
$ws.Range("N123").Value = -40225.834
# Row 132
$ws.Range("H132").Value = 3340.2258
$ws.Range("I132").Value = 2749.5
$ws.Range("J132").Value = 4158.154
$ws.Range("K132").Value = 8248.5
$ws.Range("L132").Value = 12474.462
$ws.Range("M132").Value = -5718.5
$ws.Range("N132").Value = -17534.462

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1882.9025
$ws.Range("I31").Value = 1114.3235
$ws.Range("J31").Value = 5616
$ws.Range("K31").Value = 1114.3235
$ws.Range("L31").Value = 5616
$ws.Range("M31").Value = -819.3235
$ws.Range("N31").Value = -6206
# Row 34
$ws.Range("H34").Value = 1882.9025
$ws.Range("I34").Value = 1114.3235
$ws.Range("J34").Value = 5616
$ws.Range("K34").Value = 1114.3235
$ws.Range("L34").Value = 5616
$ws.Range("M34").Value = -912.3235
$ws.Range("N34").Value = -6020
# Row 58
$ws.Range("H58").Value = 2438.55
$ws.Range("I58").Value = 2141.4167
$ws.Range("J58").Value = 2884.25
$ws.Range("K58").Value = 2141.4167
$ws.Range("L58").Value = 2884.25
$ws.Range("M58").Value = -1938.4167
$ws.Range("N58").Value = -3290.25
# Row 132
$ws.Range("H132").Value = 3184.95
$ws.Range("I132").Value = 2700.0715
$ws.Range("K132").Value = 8100.2145
$ws.Range("M132").Value = -5570.2145
# Row 134
$ws.Range("H134").Value = 1789.1621
$ws.Range("I134").Value = 1506.2188
$ws.Range("J134").Value = 3600
$ws.Range("K134").Value = 4518.6564
$ws.Range("L134").Value = 10800
$ws.Range("M134").Value = -1983.6564
$ws.Range("N134").Value = -15870
# Row 136
$ws.Range("H136").Value = 2438.55
$ws.Range("I136").Value = 2141.4167
$ws.Range("J136").Value = 2884.25
$ws.Range("K136").Value = 6424.250100000001
$ws.Range("L136").Value = 8652.75
$ws.Range("M136").Value = -3874.250100000001
$ws.Range("N136").Value = -13752.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 2310.375
$ws.Range("J4").Value = 2310.375
$ws.Range("L4").Value = 6931.125
$ws.Range("N4").Value = -7155.125
# Row 5
$ws.Range("H5").Value = 2805.3
$ws.Range("I5").Value = 3383.2856
$ws.Range("K5").Value = 10149.8568
$ws.Range("M5").Value = -10037.8568
# Row 113
$ws.Range("H113").Value = 591.8525
$ws.Range("I113").Value = 503.2927
$ws.Range("J113").Value = 773.4
$ws.Range("K113").Value = 1509.8781
$ws.Range("L113").Value = 2320.2
$ws.Range("M113").Value = 660.1218999999999
$ws.Range("N113").Value = -6660.2
# Row 131
$ws.Range("H131").Value = 18869090
$ws.Range("J131").Value = 27028472
$ws.Range("L131").Value = 81085416
$ws.Range("N131").Value = -81095496
# Row 135
$ws.Range("H135").Value = 2805.3
$ws.Range("I135").Value = 3383.2856
$ws.Range("K135").Value = 30449.5704
$ws.Range("M135").Value = -27914.5704

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3388.889
$ws.Range("I80").Value = 2720
$ws.Range("J80").Value = 4225
$ws.Range("K80").Value = 2720
$ws.Range("L80").Value = 4225
$ws.Range("M80").Value = -1722
$ws.Range("N80").Value = -6221
# Row 83
$ws.Range("H83").Value = 3388.889
$ws.Range("I83").Value = 2720
$ws.Range("J83").Value = 4225
$ws.Range("K83").Value = 13600
$ws.Range("L83").Value = 21125
$ws.Range("M83").Value = -8608
$ws.Range("N83").Value = -31109

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 733.3333
$ws.Range("I46").Value = 833.3333
$ws.Range("J46").Value = 700
$ws.Range("K46").Value = 833.3333
$ws.Range("L46").Value = 700
$ws.Range("M46").Value = -645.3333
$ws.Range("N46").Value = -1076
# Row 61
$ws.Range("H61").Value = 2754.0386
$ws.Range("I61").Value = 2584.5789
$ws.Range("J61").Value = 3214
$ws.Range("K61").Value = 2584.5789
$ws.Range("L61").Value = 3214
$ws.Range("M61").Value = -2382.5789
$ws.Range("N61").Value = -3618
# Row 82
$ws.Range("H82").Value = 2878.3333
$ws.Range("I82").Value = 1667.3334
$ws.Range("J82").Value = 3483.8333
$ws.Range("K82").Value = 1667.3334
$ws.Range("L82").Value = 3483.8333
$ws.Range("M82").Value = -1306.3334
$ws.Range("N82").Value = -4205.8333
# Row 85
$ws.Range("H85").Value = 2878.3333
$ws.Range("I85").Value = 1667.3334
$ws.Range("J85").Value = 3483.8333
$ws.Range("K85").Value = 1667.3334
$ws.Range("L85").Value = 3483.8333
$ws.Range("M85").Value = -419.3334
$ws.Range("N85").Value = -5979.8333
# Row 113
$ws.Range("H113").Value = 2754.0386
$ws.Range("I113").Value = 2584.5789
$ws.Range("J113").Value = 3214
$ws.Range("K113").Value = 2584.5789
$ws.Range("L113").Value = 3214
$ws.Range("M113").Value = -414.5789
$ws.Range("N113").Value = -7554

Write-Host "Edit complete"